$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '61.955.01'
$ws.Cells.Item(2, 5).Value = '  -0.19%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.419.52'
$ws.Cells.Item(3, 5).Value = '  +0.00%  '

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'409.68"
$ws.Cells.Item(5, 5).Value = '  +0.98%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'128.69"
$ws.Cells.Item(6, 5).Value = '  -2.86%  '

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.637"
$ws.Cells.Item(7, 5).Value = '  +7.59%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.19%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +6.53%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +10.19%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'42.83"
$ws.Cells.Item(11, 5).Value = '  +2.21%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +8.04%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -0.09%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value = '3.950.78'
$ws.Cells.Item(14, 5).Value = '  -0.32%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Chainlink'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(15, 4).Value = "'21.22"
$ws.Cells.Item(15, 5).Value = '  +7.16%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +45.66%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '3.406.11'
$ws.Cells.Item(17, 5).Value = '  -6.81%  '

# Row 18
$ws.Cells.Item(18, 4).Value = "'12.27"
$ws.Cells.Item(18, 5).Value = '  +4.54%  '

# Row 19
$ws.Cells.Item(19, 4).Value = "'1.09"
$ws.Cells.Item(19, 5).Value = '  +6.95%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '61.909.12'
$ws.Cells.Item(20, 5).Value = '  -0.18%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'452.11"
$ws.Cells.Item(21, 5).Value = '  +44.47%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'92.00"
$ws.Cells.Item(22, 5).Value = '  +10.12%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.01%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'12.92"
$ws.Cells.Item(24, 5).Value = '  +0.97%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +3.12%  '

# Row 26
$ws.Cells.Item(26, 4).Value = "'33.52"
$ws.Cells.Item(26, 5).Value = '  +12.96%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'8.85"
$ws.Cells.Item(27, 5).Value = '  +9.66%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.02%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -0.11%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'7.55"
$ws.Cells.Item(30, 5).Value = '  -3.24%  '

# Row 31
$ws.Cells.Item(31, 4).Value = "'11.99"
$ws.Cells.Item(31, 5).Value = '  +5.52%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Kaspa'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(32, 4).Value = "'0.168"
$ws.Cells.Item(32, 5).Value = '  -2.83%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(33, 4).Value = "'0.114"
$ws.Cells.Item(33, 5).Value = '  -1.13%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(34, 4).Value = "'42.78"
$ws.Cells.Item(34, 5).Value = '  +0.27%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +0.00%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +2.58%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'53.43"
$ws.Cells.Item(37, 5).Value = '  +3.83%  '

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.998"
$ws.Cells.Item(38, 5).Value = '  +0.10%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.29%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +7.53%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -1.00%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -4.58%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Monero'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(43, 4).Value = "'141.37"
$ws.Cells.Item(43, 5).Value = '  +1.47%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'NEARProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(44, 4).Value = "'4.25"
$ws.Cells.Item(44, 5).Value = '  +6.96%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.04%  '

# Row 46
$ws.Cells.Item(46, 4).Value = "'2.41"
$ws.Cells.Item(46, 5).Value = '  +8.30%  '

# Row 47
$ws.Cells.Item(47, 4).Value = "'16.53"
$ws.Cells.Item(47, 5).Value = '  -0.70%  '

# Row 48
$ws.Cells.Item(48, 4).Value = "'22.52"
$ws.Cells.Item(48, 5).Value = '  +6.06%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'ThetaToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(49, 4).Value = "'2.14"
$ws.Cells.Item(49, 5).Value = '  +9.22%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(50, 4).Value = '3.759.95'
$ws.Cells.Item(50, 5).Value = '  +0.01%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Maker'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(51, 4).Value = '2.113.35'
$ws.Cells.Item(51, 5).Value = '  +0.22%  '
